# LocalAirData.xlsx maintenance edit ("jenkins.. Ashot issue" fix).
#
# Content change:
#   Row 2 ("DataProviderWithExcel_002") has FromLocation="nyc" / ToLocation="mia".
#   The destination test data is updated from "mia" to "bost" (column E, ToLocation).
#
# UI/view-state changes captured by the same save:
#   the sheet was scrolled one column to the right (top-left visible cell becomes
#   B1) and the live selection/active cell left on D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: ToLocation "mia" -> "bost" --------------------------------
$ws.Range("E2").Value = "bost"

# --- View-state edit: scroll so column B is the left-most visible column,
#     then leave the selection/active cell on D11 -------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D11").Select() | Out-Null
